$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$cellRange = $table.Cell(1, 1).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("11÷9=1, 2", $true, $false, $false, $false, $false, $true, 1, $false, "41÷6=6, 5", 2) | Out-Null

$cellRange = $table.Cell(1, 2).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("87÷3=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "99÷8=12, 3", 2) | Out-Null

$cellRange = $table.Cell(1, 3).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("14÷3=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "71÷8=8, 7", 2) | Out-Null

$cellRange = $table.Cell(1, 4).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("12÷9=1, 3", $true, $false, $false, $false, $false, $true, 1, $false, "67÷7=9, 4", 2) | Out-Null

$cellRange = $table.Cell(1, 5).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("50÷5=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "67÷3=22, 1", 2) | Out-Null

$cellRange = $table.Cell(5, 1).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("68÷6=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "15÷2=7, 1", 2) | Out-Null

$cellRange = $table.Cell(5, 2).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("18÷8=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "26÷2=13, 0", 2) | Out-Null

$cellRange = $table.Cell(5, 3).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("74÷7=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "27÷7=3, 6", 2) | Out-Null

$cellRange = $table.Cell(5, 4).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("56÷4=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "89÷7=12, 5", 2) | Out-Null

$cellRange = $table.Cell(5, 5).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("27÷8=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "71÷5=14, 1", 2) | Out-Null

$cellRange = $table.Cell(9, 1).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("92÷5=18, 2", $true, $false, $false, $false, $false, $true, 1, $false, "81÷8=10, 1", 2) | Out-Null

$cellRange = $table.Cell(9, 2).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("76÷9=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "28÷4=7, 0", 2) | Out-Null

$cellRange = $table.Cell(9, 3).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("67÷9=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "47÷8=5, 7", 2) | Out-Null

$cellRange = $table.Cell(9, 4).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("58÷9=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "31÷8=3, 7", 2) | Out-Null

$cellRange = $table.Cell(9, 5).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("45÷4=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "55÷2=27, 1", 2) | Out-Null

$cellRange = $table.Cell(13, 1).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("11÷7=1, 4", $true, $false, $false, $false, $false, $true, 1, $false, "62÷4=15, 2", 2) | Out-Null

$cellRange = $table.Cell(13, 2).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("32÷3=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "79÷7=11, 2", 2) | Out-Null

$cellRange = $table.Cell(13, 3).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("51÷4=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "16÷4=4, 0", 2) | Out-Null

$cellRange = $table.Cell(13, 4).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("56÷2=28, 0", $true, $false, $false, $false, $false, $true, 1, $false, "11÷8=1, 3", 2) | Out-Null

$cellRange = $table.Cell(13, 5).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("90÷8=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "66÷6=11, 0", 2) | Out-Null

$cellRange = $table.Cell(17, 1).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("42÷2=21, 0", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=12, 3", 2) | Out-Null

$cellRange = $table.Cell(17, 2).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("51÷6=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "11÷9=1, 2", 2) | Out-Null

$cellRange = $table.Cell(17, 3).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("79÷6=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "95÷5=19, 0", 2) | Out-Null

$cellRange = $table.Cell(17, 4).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("19÷6=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "46÷4=11, 2", 2) | Out-Null

$cellRange = $table.Cell(17, 5).Range
$cellRange.Find.ClearFormatting()
$cellRange.Find.Execute("73÷8=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "81÷3=27, 0", 2) | Out-Null
